$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current data row (row 4), shifting the old
# data row and the source row down by one.
$ws.Rows.Item(4).Insert()

# Row 1: new title, merged across A1:I1
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Tsageri Municipality"
$ws.Range("A1:I1").Merge()

# Row 4 (newly inserted): "family with disabilities Persons " values
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value = 550
$ws.Range("C4").Value = 541
$ws.Range("D4").Value = 505
$ws.Range("E4").Value = 478
$ws.Range("F4").Value = 453
$ws.Range("G4").Value = 434
$ws.Range("H4").Value = 405
$ws.Range("I4").Value = 394

# Row 5 (former data row, now shifted down): "disabilities Persons " values
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value = 613
$ws.Range("C5").Value = 597
$ws.Range("D5").Value = 557
$ws.Range("E5").Value = 531
$ws.Range("F5").Value = 499
$ws.Range("G5").Value = 472
$ws.Range("H5").Value = 444
$ws.Range("I5").Value = 431

# Row 6 is now the source row (formerly row 5); it keeps its merge A_:H_
